$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Rows 7-12 on the "Overview" sheet share the "Latest HO Xliff Generate Date"
# timestamp - bump it from 02:24:35 to 02:24:50.
$wsOverview.Range("G7:G12").Value = "2016-09-05 02:24:50"

# "zh-cn" sheet: rows 7-12 get a Priority of "ht" and the "Latest Handoff
# Datetime" moves from 02:24:30 to 02:24:45.
$wsZhCn.Range("E7:E12").Value = "ht"
$wsZhCn.Range("H7:H12").Value = "2016-09-05 02:24:45"

# "de-de" sheet: rows 7-12 get a Priority of "ht" and the "Latest Handoff
# Datetime" (shared text with Overview's column G) moves to 02:24:50.
$wsDeDe.Range("E7:E12").Value = "ht"
$wsDeDe.Range("H7:H12").Value = "2016-09-05 02:24:50"
